{"js": "// Replace the four Munson's slogan/lema/frase sentences with their revised\n// Portuguese copy (per the commit's localization fixes), using Office.js\n// search + Replace on the full sentence text. Each sentence lives in a\n// single run/w:t in the source document, so a direct search-and-replace of\n// the entire sentence keeps the surrounding run formatting untouched.\n\nconst replacements = [\n  {\n    find: \"A campanha de marketing usar\u00e1 o seguinte slogan para capturar a ess\u00eancia da marca Munson: \\\"Munson's: Pickles and Ppreserve with a Purpose\\\".\",\n    replace: \"A campanha de marketing usar\u00e1 o seguinte slogan para capturar a ess\u00eancia da marca Munson's: \\\"Munson's: Pickles and Preserves com prop\u00f3sito\\\".\"\n  },\n  {\n    find: \"A campanha de marketing usar\u00e1 o seguinte slogan para enfatizar os benef\u00edcios do produto Munson: \\\"Munson's: More than Just Pickles and Preserves\\\".\",\n    replace: \"A campanha de marketing usar\u00e1 o seguinte slogan para enfatizar os benef\u00edcios do produto da Munson: \\\"Munson's: Mais do que apenas Pickles and Preserves\\\".\"\n  },\n  {\n    find: \"A campanha de marketing usar\u00e1 o seguinte lema para inspirar a defesa do cliente de Munson: \\\"Munson's: Share the Love of Pickles and Preserves\\\".\",\n    replace: \"A campanha de marketing usar\u00e1 o seguinte lema para inspirar a defesa do cliente da Munson's: \\\"Munson's: Compartilhe o amor por Pickles and Preserves\\\".\"\n  },\n  {\n    find: \"A campanha de marketing usar\u00e1 a seguinte frase para impulsionar o teste e a compra do produto de Munson: \\\"Munson's: Find Them, Try Them, Love Them\\\".\",\n    replace: \"A campanha de marketing usar\u00e1 a seguinte frase para impulsionar a experimenta\u00e7\u00e3o e a compra do produto da Munson's: \\\"Munson's: encontre, experimente, apaixone-se\\\".\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find target sentence: \" + find);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the four Munson's slogan/lema/frase sentences with their revised\n# Portuguese copy. Each sentence is located with Word's Find (a plain\n# search, no Replacement), and the matched Range's Text is then set\n# directly. Assigning Range.Text (rather than driving Find.Execute's\n# built-in Replace) avoids the engine's smart-quote autocorrect, so the\n# straight \" and ' characters in the new copy come through unchanged.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{\n        Old = \"A campanha de marketing usar\u00e1 o seguinte slogan para capturar a ess\u00eancia da marca Munson: `\"Munson's: Pickles and Ppreserve with a Purpose`\".\"\n        New = \"A campanha de marketing usar\u00e1 o seguinte slogan para capturar a ess\u00eancia da marca Munson's: `\"Munson's: Pickles and Preserves com prop\u00f3sito`\".\"\n    },\n    @{\n        Old = \"A campanha de marketing usar\u00e1 o seguinte slogan para enfatizar os benef\u00edcios do produto Munson: `\"Munson's: More than Just Pickles and Preserves`\".\"\n        New = \"A campanha de marketing usar\u00e1 o seguinte slogan para enfatizar os benef\u00edcios do produto da Munson: `\"Munson's: Mais do que apenas Pickles and Preserves`\".\"\n    },\n    @{\n        Old = \"A campanha de marketing usar\u00e1 o seguinte lema para inspirar a defesa do cliente de Munson: `\"Munson's: Share the Love of Pickles and Preserves`\".\"\n        New = \"A campanha de marketing usar\u00e1 o seguinte lema para inspirar a defesa do cliente da Munson's: `\"Munson's: Compartilhe o amor por Pickles and Preserves`\".\"\n    },\n    @{\n        Old = \"A campanha de marketing usar\u00e1 a seguinte frase para impulsionar o teste e a compra do produto de Munson: `\"Munson's: Find Them, Try Them, Love Them`\".\"\n        New = \"A campanha de marketing usar\u00e1 a seguinte frase para impulsionar a experimenta\u00e7\u00e3o e a compra do produto da Munson's: `\"Munson's: encontre, experimente, apaixone-se`\".\"\n    }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $found = $range.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false)\n    if (-not $found) {\n        throw \"Could not find target sentence: $($pair.Old)\"\n    }\n    $range.Text = $pair.New\n}\n"}
